$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @("order-23", "23423429", "car", "25000", "Eric", "Jamie"),
    @("23",       "24934535", "box", "12000", "tom",  "kit")
)

$startRow = 7
for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $newData[$i]
    for ($col = 1; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($rowIndex, $col)
        # Force text interpretation for numeric-looking values so they
        # are stored as text (matching the rest of the sheet), then
        # restore the default "Normal" style so no extra formatting
        # (number format / style index) is left behind on the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $rowValues[$col - 1]
        $cell.Style = "Normal"
    }
}
